$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the "sum" header's formatting (bold, border, centered) onto the
# new "Save" header cell so H1 matches the style used by the other
# header cells (G1, F1, ...).
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# New "Save" data column values.
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 1
$ws.Range("H4").Value = 0
